$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 606.4231
$ws.Range("I28").Value = 642.5714
$ws.Range("K28").Value = 642.5714
$ws.Range("M28").Value = -157.5714
$ws.Range("H33").Value = 305.42856
$ws.Range("I33").Value = 290.46155
$ws.Range("K33").Value = 290.46155
$ws.Range("M33").Value = -61.46154999999999
$ws.Range("H74").Value = 10005170
$ws.Range("I74").Value = 12503988
$ws.Range("K74").Value = 12503988
$ws.Range("M74").Value = -12503052
$ws.Range("H77").Value = 10005170
$ws.Range("I77").Value = 12503988
$ws.Range("K77").Value = 62519940
$ws.Range("M77").Value = -62515260
$ws.Range("H100").Value = 14287301
$ws.Range("I100").Value = 15385940
$ws.Range("K100").Value = 15385940
$ws.Range("M100").Value = -15385399
$ws.Range("H106").Value = 4881.364
$ws.Range("I106").Value = 2482.5
$ws.Range("K106").Value = 2482.5
$ws.Range("M106").Value = -1851.5
$ws.Range("H110").Value = 40401.332
$ws.Range("J110").Value = 40401.332
$ws.Range("L110").Value = 40401.332
$ws.Range("N110").Value = -48581.332
$ws.Range("H112").Value = 1200
$ws.Range("I112").Value = 591.75
$ws.Range("K112").Value = 1775.25
$ws.Range("M112").Value = -667.25
$ws.Range("H113").Value = 4299.2
$ws.Range("J113").Value = 4730.846
$ws.Range("L113").Value = 4730.846
$ws.Range("N113").Value = -11238.846
$ws.Range("H138").Value = 3902.02
$ws.Range("I138").Value = 803.76
$ws.Range("J138").Value = 4934.7734
$ws.Range("K138").Value = 2411.28
$ws.Range("L138").Value = 14804.3202
$ws.Range("M138").Value = 2728.72
$ws.Range("N138").Value = -25084.3202

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4245.86
$ws.Range("I32").Value = 4172.479
$ws.Range("K32").Value = 4172.479
$ws.Range("M32").Value = -3885.479
$ws.Range("H88").Value = 9527238
$ws.Range("I88").Value = 13336833
$ws.Range("K88").Value = 13336833
$ws.Range("M88").Value = -13336427
$ws.Range("H91").Value = 9527238
$ws.Range("I91").Value = 13336833
$ws.Range("K91").Value = 13336833
$ws.Range("M91").Value = -13335429
$ws.Range("H110").Value = 1664.7222
$ws.Range("I110").Value = 1519.4166
$ws.Range("K110").Value = 1519.4166
$ws.Range("M110").Value = 525.5834
$ws.Range("H132").Value = 3642.3572
$ws.Range("I132").Value = 2316
$ws.Range("J132").Value = 4379.222
$ws.Range("K132").Value = 6948
$ws.Range("L132").Value = 13137.666
$ws.Range("M132").Value = -4418
$ws.Range("N132").Value = -18197.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 835.43475
$ws.Range("I94").Value = 750.75
$ws.Range("J94").Value = 1400
$ws.Range("K94").Value = 750.75
$ws.Range("L94").Value = 1400
$ws.Range("M94").Value = -299.75
$ws.Range("N94").Value = -2302
$ws.Range("H134").Value = 2555.0264
$ws.Range("I134").Value = 2075.0435
$ws.Range("J134").Value = 3291
$ws.Range("K134").Value = 6225.130500000001
$ws.Range("L134").Value = 9873
$ws.Range("M134").Value = -3690.130500000001
$ws.Range("N134").Value = -14943

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4831906.5
$ws.Range("I16").Value = 6536843.5
$ws.Range("J16").Value = 1250.6666
$ws.Range("K16").Value = 6536843.5
$ws.Range("L16").Value = 1250.6666
$ws.Range("M16").Value = -6536556.5
$ws.Range("N16").Value = -1824.6666
$ws.Range("H99").Value = 20005902
$ws.Range("I99").Value = 66669500
$ws.Range("J99").Value = 7216.2856
$ws.Range("K99").Value = 66669500
$ws.Range("L99").Value = 7216.2856
$ws.Range("M99").Value = -66668002
$ws.Range("N99").Value = -10212.2856
$ws.Range("H113").Value = 4831906.5
$ws.Range("I113").Value = 6536843.5
$ws.Range("J113").Value = 1250.6666
$ws.Range("K113").Value = 6536843.5
$ws.Range("L113").Value = 1250.6666
$ws.Range("M113").Value = -6534673.5
$ws.Range("N113").Value = -5590.6666
$ws.Range("H126").Value = 20005902
$ws.Range("I126").Value = 66669500
$ws.Range("J126").Value = 7216.2856
$ws.Range("K126").Value = 200008500
$ws.Range("L126").Value = 21648.8568
$ws.Range("M126").Value = -200006030
$ws.Range("N126").Value = -26588.8568
$ws.Range("H134").Value = 6182.913
$ws.Range("I134").Value = 7373.6875
$ws.Range("J134").Value = 3461.1428
$ws.Range("K134").Value = 22121.0625
$ws.Range("L134").Value = 10383.4284
$ws.Range("M134").Value = -19586.0625
$ws.Range("N134").Value = -15453.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 9999.25
$ws.Range("J80").Value = 9999.25
$ws.Range("L80").Value = 29997.75
$ws.Range("N80").Value = -31869.75
$ws.Range("H83").Value = 9999.25
$ws.Range("J83").Value = 9999.25
$ws.Range("L83").Value = 89993.25
$ws.Range("N83").Value = -99353.25
$ws.Range("H92").Value = 294.16666
$ws.Range("I92").Value = 294.16666
$ws.Range("K92").Value = 882.4999799999999
$ws.Range("M92").Value = 365.5000200000001
$ws.Range("H113").Value = 530.7619
$ws.Range("I113").Value = 510.35715
$ws.Range("J113").Value = 571.5714
$ws.Range("K113").Value = 1531.07145
$ws.Range("L113").Value = 1714.7142
$ws.Range("M113").Value = 638.9285500000001
$ws.Range("N113").Value = -6054.7142
$ws.Range("H129").Value = 2401.6
$ws.Range("J129").Value = 2664.8572
$ws.Range("L129").Value = 7994.571599999999
$ws.Range("N129").Value = -17994.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3406.2727
$ws.Range("I132").Value = 3433
$ws.Range("J132").Value = 3400.3333
$ws.Range("K132").Value = 10299
$ws.Range("L132").Value = 10200.9999
$ws.Range("M132").Value = -7769
$ws.Range("N132").Value = -15260.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4495.2
$ws.Range("I7").Value = 3447.75
$ws.Range("J7").Value = 6066.375
$ws.Range("K7").Value = 3447.75
$ws.Range("L7").Value = 6066.375
$ws.Range("M7").Value = -3335.75
$ws.Range("N7").Value = -6290.375
$ws.Range("H68").Value = 649.51086
$ws.Range("I68").Value = 649.51086
$ws.Range("K68").Value = 649.51086
$ws.Range("M68").Value = 99.48914000000002
$ws.Range("H71").Value = 649.51086
$ws.Range("I71").Value = 649.51086
$ws.Range("K71").Value = 3247.5543
$ws.Range("M71").Value = 496.4457000000002
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H82").Value = 1847.0741
$ws.Range("J82").Value = 2954.2856
$ws.Range("L82").Value = 2954.2856
$ws.Range("N82").Value = -3676.2856
$ws.Range("H85").Value = 1847.0741
$ws.Range("J85").Value = 2954.2856
$ws.Range("L85").Value = 2954.2856
$ws.Range("N85").Value = -5450.2856
$ws.Range("H93").Value = 7409998
$ws.Range("I93").Value = 15874996
$ws.Range("J93").Value = 3124.875
$ws.Range("K93").Value = 15874996
$ws.Range("L93").Value = 3124.875
$ws.Range("M93").Value = -15873748
$ws.Range("N93").Value = -5620.875
$ws.Range("H126").Value = 4495.2
$ws.Range("I126").Value = 3447.75
$ws.Range("J126").Value = 6066.375
$ws.Range("K126").Value = 10343.25
$ws.Range("L126").Value = 18199.125
$ws.Range("M126").Value = -7873.25
$ws.Range("N126").Value = -23139.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20452.2
$ws.Range("J54").Value = 20452.2
$ws.Range("L54").Value = 20452.2
$ws.Range("N54").Value = -21492.2
$ws.Range("H126").Value = 2075.0476
$ws.Range("I126").Value = 1179.1818
$ws.Range("J126").Value = 3060.5
$ws.Range("K126").Value = 3537.5454
$ws.Range("L126").Value = 9181.5
$ws.Range("M126").Value = -1067.5454
$ws.Range("N126").Value = -14121.5
$ws.Range("H132").Value = 11496136
$ws.Range("I132").Value = 1223.96
$ws.Range("K132").Value = 3671.88
$ws.Range("M132").Value = -1141.88
$ws.Range("H136").Value = 1703.68
$ws.Range("I136").Value = 1432
$ws.Range("J136").Value = 2790.4
$ws.Range("K136").Value = 4296
$ws.Range("L136").Value = 8371.200000000001
$ws.Range("M136").Value = -1746
$ws.Range("N136").Value = -13471.2
